$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 already has the "Type" (A) and "Enhed" (B) values filled in, with
# bold / left-aligned formatting. Fill the same values + formatting into
# rows 3-11 of columns A and B.
$src = $ws.Range("A2:B2")
for ($r = 3; $r -le 11; $r++) {
    $addr = "A" + $r + ":B" + $r
    $dst = $ws.Range($addr)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = "Gennemsnit (nominelle priser)"
    $ws.Cells.Item($r, 2).Value = "B.1. Indestående i pengeinstitutter"
}
$excel.CutCopyMode = $false

# Move the active selection to B12, matching the post-edit view state
$ws.Range("B12").Select()
